$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price (column D) cells to remain text so values like "239.70"
# keep their exact decimal formatting instead of being parsed as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "43.922.15"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.325.37"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "97.05"
$ws.Range("E5").Value = "  +4.38%  "
$ws.Range("D6").Value = "271.09"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "45.85"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").Value = "8.11"
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "2.675.14"
$ws.Range("E14").Value = "  +4.21%  "
$ws.Range("D15").Value = "15.55"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "0.867"
$ws.Range("E16").Value = "  +8.03%  "
$ws.Range("D17").Value = "2.331.60"
$ws.Range("E17").Value = "  +4.13%  "
$ws.Range("D18").Value = "43.898.37"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  +6.67%  "
$ws.Range("D21").Value = "72.78"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").Value = "239.70"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "9.43"
$ws.Range("E24").Value = "  +4.91%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").Value = "11.37"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "3.44"
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "38.39"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "22.42"
$ws.Range("E31").Value = "  +7.45%  "
$ws.Range("D32").Value = "174.84"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "0.0907"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").Value = "0.0362"
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "4.42"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").Value = "3.37"
$ws.Range("E39").Value = "  -6.38%  "
$ws.Range("D40").Value = "0.242"
$ws.Range("E40").Value = "  +10.14%  "
$ws.Range("D41").Value = "2.36"
$ws.Range("E41").Value = "  +8.20%  "
$ws.Range("D42").Value = "1.38"
$ws.Range("E42").Value = "  +18.78%  "
$ws.Range("D43").Value = "12.23"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "9.18"
$ws.Range("E44").Value = "  +9.60%  "
$ws.Range("D45").Value = "62.27"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "5.40"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  +4.52%  "
$ws.Range("D48").Value = "100.32"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "1.20"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "0.189"
$ws.Range("E50").Value = "  +16.05%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.553.68"
$ws.Range("E51").Value = "  +4.19%  "
